$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-22 Sunday" "2024-12-23 Monday"
Replace-Text "17÷9=" "42÷5="
Replace-Text "10÷9=" "83÷8="
Replace-Text "18÷3=" "20÷6="
Replace-Text "22÷9=" "22÷6="
Replace-Text "16÷3=" "85÷9="
Replace-Text "29÷3=" "77÷6="
Replace-Text "96÷7=" "39÷9="
Replace-Text "86÷2=" "88÷4="
Replace-Text "55÷7=" "50÷4="
Replace-Text "61÷6=" "59÷4="
Replace-Text "45÷9=" "51÷6="
Replace-Text "36÷8=" "51÷5="
Replace-Text "47÷6=" "10÷2="
Replace-Text "60÷2=" "80÷8="
Replace-Text "81÷6=" "56÷9="
Replace-Text "88÷5=" "64÷9="
Replace-Text "62÷4=" "33÷5="
Replace-Text "80÷6=" "81÷5="
Replace-Text "36÷3=" "21÷2="
Replace-Text "48÷9=" "70÷5="
Replace-Text "41÷9=" "84÷5="
Replace-Text "24÷7=" "48÷7="
Replace-Text "21÷5=" "47÷7="
Replace-Text "87÷7=" "18÷9="
Replace-Text "48÷3=" "65÷4="

Write-Output "Done"
